# Apply the gh-pages data refresh (commit 456a3b4) to 合肥-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Worksheet 1): update "想去人数" (F column) counts ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7680
$ws1.Range("F3").Value = 297
$ws1.Range("F4").Value = 34
$ws1.Range("F6").Value = 4378
$ws1.Range("F7").Value = 331
$ws1.Range("F8").Value = 616
$ws1.Range("F10").Value = 690
$ws1.Range("F11").Value = 168

# ---- Sheet "演出" (Worksheet 2): ticket became unavailable + count bump ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "不可售"
$ws2.Range("F3").Value = 14

# ---- Sheet "全部类型" (Worksheet 4): combined view of all the above ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7680
$ws4.Range("G3").Value = "不可售"
$ws4.Range("F4").Value = 297
$ws4.Range("F5").Value = 34
$ws4.Range("F7").Value = 4378
$ws4.Range("F8").Value = 331
$ws4.Range("F9").Value = 616
$ws4.Range("F11").Value = 690
$ws4.Range("F12").Value = 14
$ws4.Range("F13").Value = 168
